# Update "想去人数" (number of people interested) values on the
# "展览" and "全部类型" sheets to match the freshly generated data.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 177
$ws1.Range("F4").Value  = 135
$ws1.Range("F6").Value  = 17990
$ws1.Range("F9").Value  = 1066
$ws1.Range("F10").Value = 6783
$ws1.Range("F17").Value = 150
$ws1.Range("F19").Value = 211
$ws1.Range("F25").Value = 269
$ws1.Range("F26").Value = 978
$ws1.Range("F28").Value = 5159
$ws1.Range("F30").Value = 18
$ws1.Range("F33").Value = 12008
$ws1.Range("F36").Value = 202
$ws1.Range("F37").Value = 271
$ws1.Range("F38").Value = 3909

# --- Sheet "全部类型" ---
$ws2 = $wb.Worksheets.Item("全部类型")
$ws2.Range("F2").Value  = 177
$ws2.Range("F4").Value  = 135
$ws2.Range("F6").Value  = 17990
$ws2.Range("F9").Value  = 1066
$ws2.Range("F10").Value = 6783
$ws2.Range("F17").Value = 150
$ws2.Range("F19").Value = 211
$ws2.Range("F25").Value = 269
$ws2.Range("F26").Value = 978
$ws2.Range("F28").Value = 5159
$ws2.Range("F32").Value = 18
$ws2.Range("F35").Value = 12009
$ws2.Range("F38").Value = 202
$ws2.Range("F39").Value = 271
$ws2.Range("F40").Value = 3909

$wb.Save()
